# Weekly update: a new week's price record is prepended for Albahaca
# (Primera / Segunda, Region Metropolitana, $/docena de matas) and the
# existing historical rows (492:592) are pushed down by two rows
# (to 494:594), preserving their data and formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 492:592 down by two rows (to 494:594),
# inserting two new blank rows at 492:493. Using EntireRow.Insert()
# preserves cell formatting (e.g. the date style on column D) for the
# newly created rows.
$ws.Range("A492:R493").EntireRow.Insert()

# New row for "Primera" quality, week of 2023-01-02 (serial 44932)
$ws.Range("A492").Value2 = 6
$ws.Range("B492").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C492").Value2 = "Metropolitana"
$ws.Range("D492").Value2 = 44932
$ws.Range("E492").Value2 = 13
$ws.Range("F492").Value2 = 100112052
$ws.Range("G492").Value2 = "Albahaca"
$ws.Range("H492").Value2 = "Sin especificar"
$ws.Range("I492").Value2 = "Primera"
$ws.Range("J492").Value2 = 1170
$ws.Range("K492").Value2 = 3000
$ws.Range("L492").Value2 = 4000
$ws.Range("M492").Value2 = 3530
$ws.Range("N492").Value2 = "`$/docena de matas"
$ws.Range("O492").Value2 = "Región Metropolitana"
$ws.Range("P492").Value2 = 588
$ws.Range("Q492").Value2 = 6
$ws.Range("R492").Value2 = "Hortaliza"

# New row for "Segunda" quality, same week
$ws.Range("A493").Value2 = 6
$ws.Range("B493").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C493").Value2 = "Metropolitana"
$ws.Range("D493").Value2 = 44932
$ws.Range("E493").Value2 = 13
$ws.Range("F493").Value2 = 100112052
$ws.Range("G493").Value2 = "Albahaca"
$ws.Range("H493").Value2 = "Sin especificar"
$ws.Range("I493").Value2 = "Segunda"
$ws.Range("J493").Value2 = 420
$ws.Range("K493").Value2 = 2500
$ws.Range("L493").Value2 = 3000
$ws.Range("M493").Value2 = 2810
$ws.Range("N493").Value2 = "`$/docena de matas"
$ws.Range("O493").Value2 = "Región Metropolitana"
$ws.Range("P493").Value2 = 468
$ws.Range("Q493").Value2 = 6
$ws.Range("R493").Value2 = "Hortaliza"
